# Portfolio data refresh: bump the "modified time" stamp (修改时间) on every
# row from 202509241007 -> 202509241121 across all three sheets, and update
# the 中国外运 allocation on the 大成投资组合 sheet from 32.01 -> 20.
#
# The timestamp columns are stored as TEXT (not numbers) in the workbook, so
# a plain `Range.Value = "202509241121"` assignment would get auto-coerced
# to a numeric cell by Excel's smart-entry logic. To preserve the original
# string cell type (and keep the default/general cell style untouched), we
# stage the new text through a literal-string formula and then collapse it
# to a static value via Copy + PasteSpecial(xlPasteValues) - mirroring how
# Excel itself converts a formula result to a plain cached value.

$xlPasteValues = -4163

function Set-TextValue($Sheet, $RangeAddress, $Text) {
    $range = $Sheet.Range($RangeAddress)
    $escaped = $Text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 大智投资组合 (rows 2-9, column E = 修改时间) ---
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1 "E2:E9" "202509241121"

# --- Sheet 2: 大成投资组合 (rows 2-12, column E = 修改时间) ---
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2 "E2:E12" "202509241121"
# 中国外运 (row 9) allocation percentage changed 32.01 -> 20
$ws2.Range("D9").Value = 20

# --- Sheet 3: 我的投资组合 (rows 2-13, column G = 修改时间) ---
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3 "G2:G13" "202509241121"

$excel.CutCopyMode = 0
